# Add a new "Greece" market sheet, cloned from the existing "Croatia" sheet
# (same layout/styles), fill in the Greece-specific values, and make it the
# active tab - mirroring how the "Croatia" sheet itself was most likely
# produced (copy the last country sheet, rename, fill in two cells).

$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Reset Croatia's own selection to a "whole sheet" selection before it loses
# focus, matching the state an Excel tab settles into once it is no longer
# the active sheet.
$croatia.Activate()
$croatia.Cells.Select()

# Clone Croatia (keeps all formatting, merged cells, column widths, etc.)
# and place the copy right after it - this becomes the new last sheet.
$croatia.Copy([System.Reflection.Missing]::Value, $croatia)

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Make the new sheet the active tab.
$greece.Activate()

# Fill in the Greece-specific market name and ticket reference, same cells
# that hold "Croatia Market" / "NGC-3139/T2482" on the template sheet.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3203/3186"

# Leave the selection on B4, matching where the editor's cursor ended up.
$greece.Range("B4").Select()
